$wb = $excel.ActiveWorkbook

$wsConvos = $wb.Worksheets.Item("Convos")
$wsPConvos = $wb.Worksheets.Item("PConvos")

# Update the INCLUDE sort field references from PConvos-B* to PConvos-A*
$wsConvos.Range("A2").Value = "Login please`nINCLUDE PConvos-A2"
$wsConvos.Range("A4").Value = "Logout please!`nINCLUDE PConvos-A5"

# Update the saved cell selection on each sheet, keeping "Convos" as the active tab
$wsPConvos.Range("A5").Select() | Out-Null
$wsConvos.Activate()
$wsConvos.Range("A4").Select() | Out-Null
